$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "1.00", "62.853.94") are not re-interpreted as numbers/dates.
    $r.NumberFormat = "@"
    $r.Value = $val
    # Restore the default ("Normal") style so we do not leave a stray
    # text-number-format style behind on cells that originally had none.
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "62.853.94"
Set-TextValue $ws "E2" "  +2.78%  "
Set-TextValue $ws "D3" "2.959.00"
Set-TextValue $ws "E3" "  +1.01%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.12%  "
Set-TextValue $ws "D5" "594.53"
Set-TextValue $ws "E5" "  +0.41%  "
Set-TextValue $ws "D6" "145.79"
Set-TextValue $ws "E6" "  +1.04%  "
Set-TextValue $ws "D7" "1.00"
Set-TextValue $ws "E7" "  +0.00%  "
Set-TextValue $ws "D8" "2.959.07"
Set-TextValue $ws "E8" "  +1.06%  "
Set-TextValue $ws "D9" "0.507"
Set-TextValue $ws "E9" "  +0.74%  "
Set-TextValue $ws "D10" "7.22"
Set-TextValue $ws "E10" "  +2.95%  "
Set-TextValue $ws "E11" "  +6.91%  "
Set-TextValue $ws "D12" "0.445"
Set-TextValue $ws "E12" "  +0.69%  "
Set-TextValue $ws "E13" "  +6.56%  "
Set-TextValue $ws "D14" "33.13"
Set-TextValue $ws "E14" "  -1.37%  "
Set-TextValue $ws "E15" "  -0.43%  "
Set-TextValue $ws "D16" "3.454.83"
Set-TextValue $ws "E16" "  +1.17%  "
Set-TextValue $ws "D17" "62.824.45"
Set-TextValue $ws "E17" "  +2.86%  "
Set-TextValue $ws "D18" "6.73"
Set-TextValue $ws "E18" "  +0.05%  "
Set-TextValue $ws "D19" "2.968.25"
Set-TextValue $ws "E19" "  +1.35%  "
Set-TextValue $ws "D20" "442.25"
Set-TextValue $ws "E20" "  +2.30%  "
Set-TextValue $ws "D21" "13.54"
Set-TextValue $ws "E21" "  +0.28%  "
Set-TextValue $ws "D22" "0.670"
Set-TextValue $ws "E22" "  -1.45%  "
Set-TextValue $ws "E23" "  +0.03%  "
Set-TextValue $ws "D24" "11.31"
Set-TextValue $ws "E24" "  +2.28%  "
Set-TextValue $ws "D25" "81.59"
Set-TextValue $ws "E25" "  -0.42%  "
Set-TextValue $ws "D26" "2.13"
Set-TextValue $ws "E26" "  -3.58%  "
Set-TextValue $ws "D27" "11.84"
Set-TextValue $ws "E27" "  +0.33%  "
Set-TextValue $ws "E28" "  +0.03%  "
Set-TextValue $ws "D29" "7.27"
Set-TextValue $ws "E29" "  +4.15%  "
Set-TextValue $ws "E30" "  +0.54%  "
Set-TextValue $ws "D31" "2.15"
Set-TextValue $ws "E31" "  -3.30%  "
Set-TextValue $ws "E32" "  +10.58%  "
Set-TextValue $ws "E33" "  -1.01%  "
Set-TextValue $ws "D34" "26.54"
Set-TextValue $ws "E34" "  -0.80%  "
Set-TextValue $ws "D35" "0.999"
Set-TextValue $ws "E35" "  -0.07%  "
Set-TextValue $ws "D36" "0.996"
Set-TextValue $ws "E36" "  -1.33%  "
Set-TextValue $ws "D37" "5.66"
Set-TextValue $ws "E37" "  +0.36%  "
Set-TextValue $ws "D38" "3.10"
Set-TextValue $ws "E38" "  +4.23%  "
Set-TextValue $ws "D39" "2.04"
Set-TextValue $ws "E39" "  +2.40%  "
Set-TextValue $ws "D40" "49.53"
Set-TextValue $ws "E40" "  -0.99%  "
Set-TextValue $ws "D41" "8.52"
Set-TextValue $ws "E41" "  -1.19%  "
Set-TextValue $ws "E42" "  -4.42%  "
Set-TextValue $ws "D43" "40.86"
Set-TextValue $ws "E43" "  -1.65%  "
Set-TextValue $ws "D44" "0.281"
Set-TextValue $ws "E44" "  -0.13%  "
Set-TextValue $ws "D45" "2.739.30"
Set-TextValue $ws "E45" "  +1.45%  "
Set-TextValue $ws "D46" "134.74"
Set-TextValue $ws "E46" "  +1.07%  "
Set-TextValue $ws "D47" "366.50"
Set-TextValue $ws "E47" "  -1.74%  "
Set-TextValue $ws "E48" "  -2.37%  "
Set-TextValue $ws "E49" "  +0.04%  "
Set-TextValue $ws "E50" "  -0.31%  "
Set-TextValue $ws "D51" "22.96"
Set-TextValue $ws "E51" "  -3.65%  "
